$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in B26/C26 (time-of-day entries) and let the existing shared formula
# in D26 (t="shared" si="0", which computes C-B) recompute naturally.
$ws.Range("B26").Value = 0.60416666666666663
$ws.Range("C26").Value = 0.70138888888888884

# Match the [h]:mm number format already used by the other rows in B:D
# (copy the format from a sibling cell so style indices line up).
$ws.Range("B26:C26").NumberFormat = $ws.Range("B7").NumberFormat

# New note text in E26, referencing the newly appended shared string.
$ws.Range("E26").Value = "submodules are now completely sorted out"

# Move the active selection from C26 to C27, matching the edited sheetView.
$ws.Range("C27").Select()
